# Add the newest "第613回" (No. 613) meeting row to the top of the data
# table (row 2, right below the header row), pushing every existing
# meeting row down by one — matching how a new row is published at the
# top of this "開催案内" listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at row 2. All of the existing rows
# (previously 2..44) shift down to 3..45 automatically, carrying their
# original cell contents with them unmodified.
$ws.Rows.Item(2).Insert()

# Populate the new row with the 第613回 meeting details.
$ws.Range("A2").Value = "第613回"
$ws.Range("B2").Value = "2025年7月23日（令和7年7月23日）"
$ws.Range("C2").Value = "１主な施設基準の届出状況等について`n２医療ＤＸ推進体制整備加算等の要件の見直しについて`n３入院について（その１）`n"
$ws.Range("D2").Value = "－"
$ws.Range("E2").Value = "資料`n"
$ws.Range("F2").Value = "－"

# Writing the multi-line agenda text auto-expands the row height with an
# explicit/custom height; AutoFit puts the row back to the sheet's normal
# (non-custom) auto height, matching the other rows in this table.
$ws.Rows.Item(2).AutoFit()
